$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or
# "63.710.74" are not auto-converted to numbers by Excel's smart-entry.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.710.74"
$ws.Range("E2").Value = "  +4.44%  "
$ws.Range("D3").Value = "2.727.95"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "578.15"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "154.11"
$ws.Range("E6").Value = "  +5.78%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").Value = "0.610"
$ws.Range("E8").Value = "  +1.71%  "
$ws.Range("D9").Value = "2.754.29"
$ws.Range("E9").Value = "  +2.78%  "
$ws.Range("D10").Value = "6.69"
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "0.112"
$ws.Range("E11").Value = "  +4.58%  "
$ws.Range("D12").Value = "0.161"
$ws.Range("E12").Value = "  +4.23%  "
$ws.Range("D13").Value = "0.388"
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("D14").Value = "3.221.23"
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").Value = "26.34"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("D16").Value = "63.696.50"
$ws.Range("E16").Value = "  +4.45%  "
$ws.Range("D17").Value = "0.0000152"
$ws.Range("E17").Value = "  +5.66%  "
$ws.Range("D18").Value = "2.756.94"
$ws.Range("E18").Value = "  +3.66%  "
$ws.Range("D19").Value = "11.93"
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "4.86"
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "360.50"
$ws.Range("E21").Value = "  +2.61%  "
$ws.Range("D22").Value = "6.96"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "0.532"
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "66.02"
$ws.Range("E25").Value = "  +2.95%  "
$ws.Range("E26").Value = "  +4.58%  "
$ws.Range("D27").Value = "8.51"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "0.0₃0910"
$ws.Range("E29").Value = "  +11.08%  "
$ws.Range("D30").Value = "1.99"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "7.12"
$ws.Range("E31").Value = "  +5.26%  "
$ws.Range("D32").Value = "171.42"
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("E33").Value = "  +14.35%  "
$ws.Range("E34").Value = "  -0.20%  "
$ws.Range("D35").Value = "20.44"
$ws.Range("E35").Value = "  +2.51%  "
$ws.Range("D36").Value = "4.76"
$ws.Range("E36").Value = "  +6.64%  "
$ws.Range("E37").Value = "  +8.34%  "
$ws.Range("E38").Value = "  +9.32%  "
$ws.Range("E39").Value = "  +14.32%  "
$ws.Range("D40").Value = "344.81"
$ws.Range("E40").Value = "  +4.81%  "
$ws.Range("D41").Value = "4.22"
$ws.Range("E41").Value = "  +4.93%  "
$ws.Range("D42").Value = "39.35"
$ws.Range("E42").Value = "  +2.33%  "
$ws.Range("D43").Value = "5.60"
$ws.Range("E43").Value = "  +7.24%  "
$ws.Range("D44").Value = "21.71"
$ws.Range("E44").Value = "  +5.41%  "
$ws.Range("D45").Value = "21.73"
$ws.Range("E45").Value = "  +5.58%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "0.0590"
$ws.Range("E46").Value = "  +5.00%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "139.63"
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("E48").Value = "  +4.67%  "
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "0.996"
$ws.Range("E51").Value = "  -0.13%  "
